$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh. Cell text is assigned directly; for the
# handful of "Price" values that would otherwise be auto-parsed as a
# number by Excel (e.g. "216.12", "1.00"), the cell is switched to the
# Text number format first so the original text layout is preserved.
$ws.Range('D2').Value = '27.041.28'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '1.673.57'
$ws.Range('E3').Value = '  +2.82%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.12'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.252'
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.15'
$ws.Range('E10').Value = '  +4.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +4.76%  '
$ws.Range('D12').Value = '1.909.52'
$ws.Range('E12').Value = '  +2.75%  '
$ws.Range('D13').Value = '1.673.06'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.83'
$ws.Range('E15').Value = '  +2.75%  '
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '27.057.46'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.35'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').Value = '0.0₃0737'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +3.41%  '
$ws.Range('E23').Value = '  +2.14%  '
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.41'
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.91'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0498'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').Value = '1.453.94'
$ws.Range('E33').Value = '  -4.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('E34').Value = '  +5.10%  '
$ws.Range('E35').Value = '  +5.33%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.570'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('E38').Value = '  +6.88%  '
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +3.54%  '
$ws.Range('E42').Value = '  +11.91%  '
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.87'
$ws.Range('E44').Value = '  +4.30%  '
$ws.Range('D45').Value = '1.819.12'
$ws.Range('E45').Value = '  +2.82%  '
$ws.Range('E46').Value = '  +2.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.30'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('E49').Value = '  +4.17%  '
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('E51').Value = '  +1.91%  '
